$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "슈퍼짱짱"

$ws.Range("D9").Value = "서울소재 데이터 사이언스 대학원 다니는 분 수업과제 현실"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/ds-grad-school-in-class-exer/#utm_source=rss&utm_medium=rss&utm_campaign=ds-grad-school-in-class-exer"

$ws.Range("D19").Value = "아기 있는 집에는 에몬스홈 그란데 가죽소파 4인"

$ws.Range("D44").Value = "Engineer-Ladder"

$ws.Range("D50").Value = "파이썬 3.11.0"
$ws.Range("E50").Value = "http://incredible.egloos.com/7572513"

$ws.Range("D51").Value = "[python] 문자열 대소문자 변환하기 (upper, lower, isupper, islower)"
$ws.Range("E51").Value = "https://bskyvision.com/entry/python-%EB%AC%B8%EC%9E%90%EC%97%B4-%EB%8C%80%EC%86%8C%EB%AC%B8%EC%9E%90-%EB%B3%80%ED%99%98%ED%95%98%EA%B8%B0-upper-lower-isupper-islower"
